$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume-change (E) columns are stored as plain text in
# this workbook (thousands separators written with dots, fixed decimal
# places, and "  +x.xx%  " strings padded with spaces). Writing a
# numeric-looking string straight into .Value lets Excel re-interpret it as
# a real number, silently dropping meaningful trailing zeros (e.g. "2.310"
# -> 2.31). Prefixing with an apostrophe forces text entry; resetting the
# cell's Style back to "Normal" afterwards clears the quote-prefix flag Excel
# stamps on the cell so formatting stays identical to the original file.

$ws.Range('D2').Value = '''29.647.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -3.18%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.093.73'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -1.22%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''1.009'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.23%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''342.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -2.17%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '''  -0.18%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.5118'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -2.71%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.4395'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -2.63%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''53.34'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -2.12%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.09123'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +0.52%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''1.168'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -0.88%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''24.76'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +0.70%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''2.104.83'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +0.23%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''6.731'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -1.65%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''8.177'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +1.17%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''99.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -3.02%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  -2.51%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.28%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''21.07'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +8.21%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''0.06641'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -1.14%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''1.007'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -0.23%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''6.165'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -2.48%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''29.712.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '''  -2.07%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''2.310'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -3.14%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''2.346.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.40%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''21.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -3.14%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''162.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -1.59%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''2.514'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -1.60%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''132.38'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -3.07%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''1.127'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -5.72%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D33').Value = '''1.631'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -1.86%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''6.151'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -3.69%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''3.961'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -1.33%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''6.030'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +1.81%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''10.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -1.41%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.02569'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -2.97%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.06663'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -3.02%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '''12.36'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -1.88%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.6845'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.97%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '''0.2224'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -4.17%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''1.293'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +1.48%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''0.6678'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +3.23%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''14.15'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -4.38%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''2.288'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -2.00%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''3.607'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -3.85%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '''  -2.80%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '''0.00000000335'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -8.13%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '''81.60'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -1.27%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = '''  -2.57%  '
$ws.Range('E51').Style = 'Normal'
